$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("D2").Value = "'68.566.85"
    $ws.Range("D2").Style = "Normal"
    $ws.Range("E2").Value = "'  +0.57%  "
    $ws.Range("E2").Style = "Normal"
    $ws.Range("D3").Value = "'2.699.75"
    $ws.Range("D3").Style = "Normal"
    $ws.Range("E3").Value = "'  +2.21%  "
    $ws.Range("E3").Style = "Normal"
    $ws.Range("E4").Value = "'  -0.01%  "
    $ws.Range("E4").Style = "Normal"
    $ws.Range("D5").Value = "'598.67"
    $ws.Range("D5").Style = "Normal"
    $ws.Range("E5").Value = "'  +0.26%  "
    $ws.Range("E5").Style = "Normal"
    $ws.Range("D6").Value = "'160.16"
    $ws.Range("D6").Style = "Normal"
    $ws.Range("E6").Value = "'  +2.61%  "
    $ws.Range("E6").Style = "Normal"
    $ws.Range("E7").Value = "'  +0.05%  "
    $ws.Range("E7").Style = "Normal"
    $ws.Range("D8").Value = "'0.544"
    $ws.Range("D8").Style = "Normal"
    $ws.Range("E8").Value = "'  +0.36%  "
    $ws.Range("E8").Style = "Normal"
    $ws.Range("D9").Value = "'2.698.21"
    $ws.Range("D9").Style = "Normal"
    $ws.Range("E9").Value = "'  +2.21%  "
    $ws.Range("E9").Style = "Normal"
    $ws.Range("E10").Value = "'  +0.09%  "
    $ws.Range("E10").Style = "Normal"
    $ws.Range("E11").Value = "'  -0.32%  "
    $ws.Range("E11").Style = "Normal"
    $ws.Range("D12").Value = "'5.30"
    $ws.Range("D12").Style = "Normal"
    $ws.Range("E12").Value = "'  +1.09%  "
    $ws.Range("E12").Style = "Normal"
    $ws.Range("D13").Value = "'0.360"
    $ws.Range("D13").Style = "Normal"
    $ws.Range("E13").Value = "'  +2.71%  "
    $ws.Range("E13").Style = "Normal"
    $ws.Range("D14").Value = "'28.27"
    $ws.Range("D14").Style = "Normal"
    $ws.Range("E14").Value = "'  +1.14%  "
    $ws.Range("E14").Style = "Normal"
    $ws.Range("D15").Value = "'3.190.69"
    $ws.Range("D15").Style = "Normal"
    $ws.Range("E15").Value = "'  +2.20%  "
    $ws.Range("E15").Style = "Normal"
    $ws.Range("D16").Value = "'0.0000188"
    $ws.Range("D16").Style = "Normal"
    $ws.Range("E16").Value = "'  -0.82%  "
    $ws.Range("E16").Style = "Normal"
    $ws.Range("D17").Value = "'68.639.75"
    $ws.Range("D17").Style = "Normal"
    $ws.Range("E17").Value = "'  +0.70%  "
    $ws.Range("E17").Style = "Normal"
    $ws.Range("D18").Value = "'2.689.93"
    $ws.Range("D18").Style = "Normal"
    $ws.Range("E18").Value = "'  +1.93%  "
    $ws.Range("E18").Style = "Normal"
    $ws.Range("E19").Value = "'  +4.22%  "
    $ws.Range("E19").Style = "Normal"
    $ws.Range("D20").Value = "'365.33"
    $ws.Range("D20").Style = "Normal"
    $ws.Range("E20").Value = "'  +0.57%  "
    $ws.Range("E20").Style = "Normal"
    $ws.Range("E21").Value = "'  +3.96%  "
    $ws.Range("E21").Style = "Normal"
    $ws.Range("E22").Value = "'  +2.38%  "
    $ws.Range("E22").Style = "Normal"
    $ws.Range("E23").Value = "'  +2.52%  "
    $ws.Range("E23").Style = "Normal"
    $ws.Range("E24").Value = "'  +2.08%  "
    $ws.Range("E24").Style = "Normal"
    $ws.Range("D25").Value = "'74.34"
    $ws.Range("D25").Style = "Normal"
    $ws.Range("E25").Value = "'  -1.38%  "
    $ws.Range("E25").Style = "Normal"
    $ws.Range("D27").Value = "'9.91"
    $ws.Range("D27").Style = "Normal"
    $ws.Range("E27").Value = "'  +1.80%  "
    $ws.Range("E27").Style = "Normal"
    $ws.Range("E28").Value = "'  +2.27%  "
    $ws.Range("E28").Style = "Normal"
    $ws.Range("D29").Value = "'0.0000105"
    $ws.Range("D29").Style = "Normal"
    $ws.Range("E29").Value = "'  +0.69%  "
    $ws.Range("E29").Style = "Normal"
    $ws.Range("B30").Value = "'Bittensor"
    $ws.Range("B30").Style = "Normal"
    $ws.Range("C30").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
    $ws.Range("C30").Style = "Normal"
    $ws.Range("D30").Value = "'580.85"
    $ws.Range("D30").Style = "Normal"
    $ws.Range("E30").Value = "'  +4.55%  "
    $ws.Range("E30").Style = "Normal"
    $ws.Range("B31").Value = "'Binance-PegBSC-USD"
    $ws.Range("B31").Style = "Normal"
    $ws.Range("C31").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
    $ws.Range("C31").Style = "Normal"
    $ws.Range("D31").Value = "'0.999"
    $ws.Range("D31").Style = "Normal"
    $ws.Range("E31").Value = "'  -5.58%  "
    $ws.Range("E31").Style = "Normal"
    $ws.Range("D32").Value = "'8.23"
    $ws.Range("D32").Style = "Normal"
    $ws.Range("E32").Value = "'  +2.34%  "
    $ws.Range("E32").Style = "Normal"
    $ws.Range("D33").Value = "'1.44"
    $ws.Range("D33").Style = "Normal"
    $ws.Range("E33").Value = "'  +2.40%  "
    $ws.Range("E33").Style = "Normal"
    $ws.Range("E34").Value = "'  +5.56%  "
    $ws.Range("E34").Style = "Normal"
    $ws.Range("E35").Value = "'  +3.89%  "
    $ws.Range("E35").Style = "Normal"
    $ws.Range("E36").Value = "'  +6.28%  "
    $ws.Range("E36").Style = "Normal"
    $ws.Range("D38").Value = "'161.71"
    $ws.Range("D38").Style = "Normal"
    $ws.Range("E38").Value = "'  +0.09%  "
    $ws.Range("E38").Style = "Normal"
    $ws.Range("D39").Value = "'19.80"
    $ws.Range("D39").Style = "Normal"
    $ws.Range("E39").Value = "'  +0.80%  "
    $ws.Range("E39").Style = "Normal"
    $ws.Range("E40").Value = "'  +2.19%  "
    $ws.Range("E40").Style = "Normal"
    $ws.Range("D41").Value = "'1.91"
    $ws.Range("D41").Style = "Normal"
    $ws.Range("E41").Value = "'  +2.43%  "
    $ws.Range("E41").Style = "Normal"
    $ws.Range("D42").Value = "'5.40"
    $ws.Range("D42").Style = "Normal"
    $ws.Range("E42").Value = "'  +1.48%  "
    $ws.Range("E42").Style = "Normal"
    $ws.Range("D43").Value = "'2.68"
    $ws.Range("D43").Style = "Normal"
    $ws.Range("E43").Value = "'  +2.79%  "
    $ws.Range("E43").Style = "Normal"
    $ws.Range("E44").Value = "'  +0.38%  "
    $ws.Range("E44").Style = "Normal"
    $ws.Range("D45").Value = "'0.0₆0317"
    $ws.Range("D45").Style = "Normal"
    $ws.Range("E45").Value = "'  -4.75%  "
    $ws.Range("E45").Style = "Normal"
    $ws.Range("E46").Value = "'  +0.04%  "
    $ws.Range("E46").Style = "Normal"
    $ws.Range("D47").Value = "'157.81"
    $ws.Range("D47").Style = "Normal"
    $ws.Range("E47").Value = "'  -0.45%  "
    $ws.Range("E47").Style = "Normal"
    $ws.Range("D48").Value = "'3.94"
    $ws.Range("D48").Style = "Normal"
    $ws.Range("E48").Value = "'  +5.96%  "
    $ws.Range("E48").Style = "Normal"
    $ws.Range("E49").Value = "'  +5.05%  "
    $ws.Range("E49").Style = "Normal"
    $ws.Range("E50").Value = "'  +7.15%  "
    $ws.Range("E50").Style = "Normal"
    $ws.Range("D51").Value = "'22.03"
    $ws.Range("D51").Style = "Normal"
    $ws.Range("E51").Value = "'  -0.07%  "
    $ws.Range("E51").Style = "Normal"
